$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# 1) Rename the header row: "<Feld>_old" -> "<Feld>_FV2410"
#    and "<Feld>_new" -> "<Feld>_FV2504" (the "diff" header is kept).
# ------------------------------------------------------------------
$headerMap = @{
    "A1" = "Segmentname_FV2410"
    "B1" = "Segmentgruppe_FV2410"
    "C1" = "Segment_FV2410"
    "D1" = "Datenelement_FV2410"
    "E1" = "Segment ID_FV2410"
    "F1" = "Code_FV2410"
    "G1" = "Qualifier_FV2410"
    "H1" = "Beschreibung_FV2410"
    "I1" = "Bedingungsausdruck_FV2410"
    "J1" = "Bedingung_FV2410"
    "K1" = "diff"
    "L1" = "Segmentname_FV2504"
    "M1" = "Segmentgruppe_FV2504"
    "N1" = "Segment_FV2504"
    "O1" = "Datenelement_FV2504"
    "P1" = "Segment ID_FV2504"
    "Q1" = "Code_FV2504"
    "R1" = "Qualifier_FV2504"
    "S1" = "Beschreibung_FV2504"
    "T1" = "Bedingungsausdruck_FV2504"
    "U1" = "Bedingung_FV2504"
}

foreach ($addr in $headerMap.Keys) {
    $ws.Range($addr).Value = $headerMap[$addr]
}

# ------------------------------------------------------------------
# 2) Turn the used range into an Excel table ("Table1") so the header
#    row gets filter buttons and the whole range becomes a ListObject.
# ------------------------------------------------------------------
$tbl = $ws.ListObjects.Add(1, $ws.Range("A1:U70"), $null, 1)
$tbl.Name = "Table1"
$tbl.TableStyle = ""

# Re-apply the original header formatting (bold, grey fill, borders,
# centered + wrapped) since adding the table resets the header style.
$hdr = $ws.Range("A1:U1")
$hdr.Font.Bold = $true
$hdr.Interior.Color = 14277081
$hdr.HorizontalAlignment = -4108
$hdr.WrapText = $true
$hdr.Borders.LineStyle = 1
$hdr.Borders.Weight = 2

# ------------------------------------------------------------------
# 3) Freeze the header row (split below row 1).
# ------------------------------------------------------------------
[void]$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
